$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns for touched rows to remain Text,
# matching the inline-string cell type used throughout the sheet.
$ws.Range("D2:E17").NumberFormat = "@"
$ws.Range("D19:E25").NumberFormat = "@"
$ws.Range("D38:E48").NumberFormat = "@"
$ws.Range("D50:E51").NumberFormat = "@"

# Apply the updated coin data (re-ranked rows + refreshed price/volume figures).
$ws.Range("D2").Value = "309.24"
$ws.Range("E2").Value = "0.00%"
$ws.Range("D3").Value = "41.19"
$ws.Range("E3").Value = "-0.08%"
$ws.Range("D4").Value = "5.189"
$ws.Range("E4").Value = "1.08%"
$ws.Range("D5").Value = "0.07679"
$ws.Range("E5").Value = "0.67%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "1.678"
$ws.Range("E6").Value = "3.36%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "0.9147"
$ws.Range("E7").Value = "1.00%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "2.425"
$ws.Range("E8").Value = "-2.37%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1239"
$ws.Range("E9").Value = "11.47%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1827"
$ws.Range("E10").Value = "1.86%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.09128"
$ws.Range("E11").Value = "0.96%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.04195"
$ws.Range("E12").Value = "-1.55%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.1052"
$ws.Range("E13").Value = "0.23%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001302"
$ws.Range("E14").Value = "3.41%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005760"
$ws.Range("E15").Value = "1.25%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.345"
$ws.Range("E16").Value = "0.10%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.305"
$ws.Range("E17").Value = "1.32%"
$ws.Range("D19").Value = "7.411"
$ws.Range("E19").Value = "11.04%"
$ws.Range("D20").Value = "0.1376"
$ws.Range("E20").Value = "0.91%"
$ws.Range("D21").Value = "0.2826"
$ws.Range("E21").Value = "4.40%"
$ws.Range("D22").Value = "0.04018"
$ws.Range("E22").Value = "-0.27%"
$ws.Range("D23").Value = "0.001270"
$ws.Range("E23").Value = "0.98%"
$ws.Range("D24").Value = "0.004088"
$ws.Range("E24").Value = "-0.16%"
$ws.Range("D25").Value = "0.0001304"
$ws.Range("E25").Value = "0.36%"
$ws.Range("D38").Value = "0.02524"
$ws.Range("E38").Value = "4.27%"
$ws.Range("D39").Value = "0.05309"
$ws.Range("E39").Value = "1.30%"
$ws.Range("D40").Value = "0.007841"
$ws.Range("E40").Value = "0.33%"
$ws.Range("D41").Value = "0.1307"
$ws.Range("E41").Value = "0.52%"
$ws.Range("D42").Value = "0.006658"
$ws.Range("E42").Value = "-5.52%"
$ws.Range("D43").Value = "0.001876"
$ws.Range("E43").Value = "-3.76%"
$ws.Range("D44").Value = "0.007413"
$ws.Range("E44").Value = "-12.14%"
$ws.Range("D45").Value = "0.3066"
$ws.Range("E45").Value = "-8.11%"
$ws.Range("D46").Value = "0.00006800"
$ws.Range("E46").Value = "-1.29%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.36%"
$ws.Range("D48").Value = "0.2335"
$ws.Range("E48").Value = "293.03%"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").Value = "0.36%"
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").Value = "0.36%"
